# Renumber the API chapter titles from "10.x"/"11.x" to "12.x".
#
# Each affected title lives in the ctrTitle placeholder (always
# Shapes.Item(1) on these slides). On slides 3, 7 and 8 the whole
# title is a single run; on slides 4, 5 and 6 the chapter-number
# prefix is the first run of a two-run paragraph (the second run,
# e.g. "iconlabel", must stay untouched).
#
# Replacing text through Characters(start, length) exactly at the
# existing run's character-length boundary swaps that run's text in
# place without splitting/merging runs, so formatting (color, bold,
# size, ...) on every run - including the untouched second run - is
# preserved.

$p = $ppt.ActivePresentation

$targets = @(
    @{ Slide = 3; Length = 14; Text = "12.1 API: i18n" },
    @{ Slide = 4; Length = 10; Text = "12.1 API: " },
    @{ Slide = 5; Length = 10; Text = "12.2 API: " },
    @{ Slide = 6; Length = 10; Text = "12.2 API: " },
    @{ Slide = 7; Length = 16; Text = "12.3 API: images" },
    @{ Slide = 8; Length = 16; Text = "12.3 API: images" }
)

foreach ($t in $targets) {
    $s = $p.Slides.Item($t.Slide)
    $tr = $s.Shapes.Item(1).TextFrame.TextRange
    $tr.Characters(1, $t.Length).Text = $t.Text
}
